$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 11998
$wsExhibit.Range("F6").Value = 364
$wsExhibit.Range("F8").Value = 11899
$wsExhibit.Range("F12").Value = 577
$wsExhibit.Range("F14").Value = 5889
$wsExhibit.Range("F16").Value = 3550
$wsExhibit.Range("F18").Value = 27

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 11998
$wsAll.Range("F9").Value = 364
$wsAll.Range("F11").Value = 11899
$wsAll.Range("F15").Value = 577
$wsAll.Range("F18").Value = 5889
$wsAll.Range("F20").Value = 3550
$wsAll.Range("F22").Value = 27
